$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Consumption CE (column B) and Total Wealth CE (column C) values,
# recalculated using income_test data.
$ws.Range("B2").Value = 15470.62290791229
$ws.Range("C2").Value = 847722.6056723457

$ws.Range("B3").Value = 23739.24350645736
$ws.Range("C3").Value = 1300806.921723364

$ws.Range("B4").Value = 38346.253017191
$ws.Range("C4").Value = 2101207.282925812
